# bp - SSO herschrijven + HSTS
# Insert a new row (7) with a NetScaler Management entry (IP 10.1.205.7/16),
# shifting the rows below it down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(7).Insert()

$nbsp = [char]0x00A0
$ws.Cells.Item(7,1).Value = "IBM x3650 M3" + $nbsp + "7945G2G"
$ws.Cells.Item(7,1).Style = "Normal"
$ws.Cells.Item(7,2).Value = 1
$ws.Cells.Item(7,3).Value = "10.1.205.7/16"
$ws.Cells.Item(7,4).Value = "netscaler"
$ws.Cells.Item(7,5).Value = "Intern"
$ws.Cells.Item(7,6).Value = "10.1.5.32"
$ws.Cells.Item(7,7).Value = "VIP NetScaler 12"
$ws.Cells.Item(7,8).Value = "NetScaler Management"

$ws.Range("C8").Select() | Out-Null
